# Apply odds updates to Sheet1 per the commit diff (Atualizando o arquivo XLSX)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
# Row 4
$ws.Range("AH4").Value = 17
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 41
$ws.Range("AP4").Value = 26
$ws.Range("AT4").Value = 2.38
$ws.Range("G4").Value = 2.1
$ws.Range("I4").Value = 3.8
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 9
# Row 8
$ws.Range("AI8").Value = 10
$ws.Range("AJ8").Value = 23
$ws.Range("AN8").Value = 4.75
$ws.Range("AU8").Value = 8.5
$ws.Range("AX8").Value = 15
$ws.Range("AY8").Value = 26
$ws.Range("AZ8").Value = 51
$ws.Range("G8").Value = 2.9
$ws.Range("I8").Value = 2.35
$ws.Range("J8").Value = 3.6
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("W8").Value = 8.5
$ws.Range("Y8").Value = 11
$ws.Range("Z8").Value = 29
# Row 12
$ws.Range("AG12").Value = 15
$ws.Range("AZ12").Value = 201
$ws.Range("G12").Value = 1.42
$ws.Range("H12").Value = 4.1
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 2.2
$ws.Range("L12").Value = 8
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
# Row 14
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 7
